# Iraq League base update (17-02-2024 22:47)
# The source data rows got re-sorted/re-matched, which manifests as pairs of
# adjacent rows swapping all of their match-data columns (B, and F through AC)
# while keeping the row's own id (column A) and the Div/Div Original
# Name/Date columns (C, D, E - identical between the two rows anyway) fixed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param(
        [int]$RowA,
        [int]$RowB
    )

    $rangeA = $ws.Range("B$RowA`:AC$RowA")
    $rangeB = $ws.Range("B$RowB`:AC$RowB")

    $valuesA = $rangeA.Value2
    $valuesB = $rangeB.Value2

    $rangeA.Value2 = $valuesB
    $rangeB.Value2 = $valuesA
}

Swap-RowData -RowA 14 -RowB 15
Swap-RowData -RowA 16 -RowB 17
Swap-RowData -RowA 137 -RowB 138
Swap-RowData -RowA 141 -RowB 142
Swap-RowData -RowA 161 -RowB 162
Swap-RowData -RowA 172 -RowB 173

Write-Output "Row swaps applied"
